$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

# 1) Swap the envelope-emoji prefix for a plain envelope glyph wherever it
#    shows up (the "Email Address" column for each employee row).
$replacedAny = $false
for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -like "*📧*") {
            $cell.Value = $val -replace "📧", "✉"
            $replacedAny = $true
        }
    }
}

# 2) Column C (Email Address) is a touch too wide now that the glyph is
#    narrower — bring it in by one character width.
if ($replacedAny) {
    $ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(3).ColumnWidth - 1
}

# 3) Row 5 was missing the explicit row-height the other data rows carry
#    (it was silently falling back to the sheet default) — give it the
#    same custom height as its neighbours.
for ($r = 1; $r -le $rowCount; $r++) {
    if ($ws.Rows.Item($r).RowHeight -ne 13) {
        $ws.Rows.Item($r).RowHeight = 13
    }
}
